{"js": "// 1) \"I have downloaded the Boston Crime dataset from Kaggle.\" ->\n//    \"I am using the Boston Crime dataset from Kaggle.\"\nconst body = context.document.body;\n\nconst downloadedHits = body.search(\"I have downloaded the Boston Crime dataset from Kaggle.\", { matchCase: true });\ndownloadedHits.load(\"items\");\nawait context.sync();\n\nif (downloadedHits.items.length > 0) {\n  downloadedHits.items[0].insertText(\"I am using the Boston Crime dataset from Kaggle.\", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"I will try and find out which day of the week has the highest crime rate.\" ->\n//    \"...and see if the crimes occur during the day or the night.\"\nconst dayHits = body.search(\"I will try and find out which day of the week has the highest crime rate.\", { matchCase: true });\ndayHits.load(\"items\");\nawait context.sync();\n\nlet lastListParagraph = null;\nif (dayHits.items.length > 0) {\n  const dayRange = dayHits.items[0];\n  dayRange.insertText(\n    \"I will try and find out which day of the week has the highest crime rate and see if the crimes occur during the day or the night.\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  lastListParagraph = dayRange.paragraphs.getFirst();\n  lastListParagraph.load(\"text\");\n  await context.sync();\n}\n\n// 3) Add two new bullet-list paragraphs (same list as the \"Approach\" bullets)\n//    right after the \"day of the week\" bullet.\nif (lastListParagraph) {\n  const newPara1 = lastListParagraph.insertParagraph(\n    \"Idnetify the top crimes that occur in the city and find the areas affected by them.\",\n    \"After\"\n  );\n  await context.sync();\n\n  newPara1.insertParagraph(\n    \"Conclude the project by summarizing my findings.\",\n    \"After\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"I have downloaded the Boston Crime dataset from Kaggle.\" ->\n#    \"I am using the Boston Crime dataset from Kaggle.\"\n$downloadedRange = $d.Content\n$downloadedFound = $downloadedRange.Find.Execute(\"I have downloaded the Boston Crime dataset from Kaggle.\")\nif ($downloadedFound) {\n    $downloadedRange.Text = \"I am using the Boston Crime dataset from Kaggle.\"\n}\n\n# 2) \"I will try and find out which day of the week has the highest crime rate.\" ->\n#    \"...and see if the crimes occur during the day or the night.\"\n$dayRange = $d.Content\n$dayFound = $dayRange.Find.Execute(\"I will try and find out which day of the week has the highest crime rate.\")\nif ($dayFound) {\n    $dayRange.Text = \"I will try and find out which day of the week has the highest crime rate and see if the crimes occur during the day or the night.\"\n\n    # 3) Add two new bullet-list paragraphs (same list as the \"Approach\" bullets)\n    #    right after the \"day of the week\" bullet.\n    $targetPara = $dayRange.Paragraphs.Item(1)\n    $targetStart = $targetPara.Range.Start\n\n    $count = $d.Paragraphs.Count\n    $targetIndex = -1\n    for ($i = 1; $i -le $count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Start -eq $targetStart) {\n            $targetIndex = $i\n        }\n    }\n\n    if ($targetIndex -gt 0) {\n        $targetPara.Range.InsertParagraphAfter()\n        $newPara1 = $d.Paragraphs.Item($targetIndex + 1)\n        $newPara1.Range.Text = \"Idnetify the top crimes that occur in the city and find the areas affected by them.\"\n\n        $newPara1.Range.InsertParagraphAfter()\n        $newPara2 = $d.Paragraphs.Item($targetIndex + 2)\n        $newPara2.Range.Text = \"Conclude the project by summarizing my findings.\"\n    }\n}\n"}
